# offset.xlsx edit: add "Power Current Limit" / "Power meter probe Channel"
# columns (E/F shifted old E-H into G-J) plus a new "Select Freq" column K on
# the config sheet, add matching legend rows + a new "Column" legend column
# on the help sheet, and refresh column widths / selections to match.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # "config"
$ws2 = $wb.Worksheets.Item(2)   # "help"

# ---------------------------------------------------------------------
# Sheet "config" (sheet1): shift old E:H header block right to G:J, add
# new E:F (Power Current Limit / Power meter probe Channel) header block,
# add K (Select Freq) header + per-row values, add I3/J3 + data rows.
# Cell write order matters: it controls the order brand-new strings are
# appended to the shared-string table, so it is kept deliberate below.
# ---------------------------------------------------------------------

$cfgPairs = @(
  @("K1", "Select Freq"),
  @("K2", "GHz"),
  @("E2", "A"),
  @("E1", "Power Current Limit"),
  @("F2", "Channel"),
  @("F1", "Power meter probe Channel"),
  @("G1", "Aging time"),
  @("H1", "Source generate"),
  @("I1", "Power supply"),
  @("J1", "Power meter"),
  @("G2", "sec"),
  @("I2", "GPIB"),
  @("J2", "GPIB"),
  @("E3", 28),
  @("F3", 2),
  @("H3", 9),
  @("I3", 1),
  @("J3", 13),
  @("K3", 6),
  @("K4", 6.5),
  @("K5", 7),
  @("K6", 7.5),
  @("K7", 8),
  @("K8", 8.5),
  @("K9", 9),
  @("K10", 9.5),
  @("K11", 10),
  @("K12", 10.5),
  @("K13", 11)
)
foreach ($p in $cfgPairs) {
  $ws1.Range($p[0]).Value = $p[1]
}

# ---------------------------------------------------------------------
# Sheet "help" (sheet2): insert a "Column" legend column (B) before the
# old unit column, add an "EAC" column (E) carrying "fixed" notes, and
# add two new legend rows (Power Current Limit / Power meter probe
# Channel). As above, write order is deliberate for shared-string order.
# ---------------------------------------------------------------------

$helpPairs = @(
  @("D6", "power supply current limit data"),
  @("D5", "power supply voltage data"),
  @("D7", "power meter probe number "),
  @("E3", "fixed"),
  @("E4", "fixed"),
  @("E5", "fixed"),
  @("E6", "fixed"),
  @("E7", "fixed"),
  @("B1", "Column"),
  @("B3", "B"),
  @("B4", "C"),
  @("B5", "D"),
  @("B6", "E"),
  @("B7", "F"),
  @("B8", "G"),
  @("B9", "H"),
  @("B10", "I"),
  @("B11", "J"),
  @("C1", "unit(row = 2)"),
  @("A1", "name(row = 1)"),
  @("D1", "description"),
  @("E1", "EAC"),
  @("B2", "A"),
  @("C2", "Hz, KHz, MHz, GHz, THz"),
  @("D2", "source generate freq data"),
  @("C3", "dB"),
  @("D3", "power meter input offset data"),
  @("C4", "dB"),
  @("D4", "power meter output offset data"),
  @("C5", "V"),
  @("A6", "Power Current Limit"),
  @("C6", "A"),
  @("A7", "Power meter probe Channel"),
  @("C7", "Channel"),
  @("A8", "Aging time"),
  @("C8", "sec, min, hour"),
  @("D8", "Aging at start frequency"),
  @("A9", "Source generate"),
  @("C9", "GPIB, USB, Serial"),
  @("D9", "source generate communication data"),
  @("E9", "if select serial, Baudrate data needed"),
  @("A10", "Power supply"),
  @("C10", "GPIB, USB, Serial"),
  @("D10", "power supply communication data"),
  @("E10", "if select serial, Baudrate data needed"),
  @("A11", "Power meter"),
  @("C11", "GPIB, USB, Serial"),
  @("D11", "power meter communication dataw"),
  @("E11", "if select serial, Baudrate data needed")
)
foreach ($p in $helpPairs) {
  $ws2.Range($p[0]).Value = $p[1]
}

# E1 is a brand-new header cell on row 1 - copy the shaded header format
# from an existing row-1 header cell so it reuses the same cell style.
$ws2.Range("D1").Copy()
$ws2.Range("E1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------
# Column widths (best-effort; the headless engine quantizes widths to a
# pixel grid that doesn't perfectly reproduce Excel's font-metric based
# bestFit numbers, so these are the closest achievable values).
# ---------------------------------------------------------------------

$ws1.Columns.Item(1).ColumnWidth = 11.883370535714286
$ws1.Columns.Item(2).ColumnWidth = 10.285714285714286
$ws1.Columns.Item(3).ColumnWidth = 11.883370535714286
$ws1.Columns.Item(4).ColumnWidth = 12.883370535714286
$ws1.Columns.Item(5).ColumnWidth = 17.684151785714285
$ws1.Columns.Item(6).ColumnWidth = 24.984933035714285
$ws1.Columns.Item(7).ColumnWidth = 9.883370535714286
$ws1.Columns.Item(8).ColumnWidth = 14.484933035714286
$ws1.Columns.Item(9).ColumnWidth = 11.883370535714286
$ws1.Columns.Item(10).ColumnWidth = 11.285714285714286
$ws1.Columns.Item(11).ColumnWidth = 9.642857142857144

$ws2.Columns.Item(1).ColumnWidth = 24.883370535714285
$ws2.Columns.Item(2).ColumnWidth = 6.883370535714286
$ws2.Columns.Item(3).ColumnWidth = 21.285714285714285
$ws2.Columns.Item(4).ColumnWidth = 32.984933035714285
$ws2.Columns.Item(5).ColumnWidth = 32.684151785714285

# help sheet gains an explicit page setup (matching the config sheet's).
$ws2.PageSetup.PaperSize = 9
$ws2.PageSetup.Orientation = 1

# ---------------------------------------------------------------------
# View state: config sheet stays the active tab with K4 selected; help
# sheet keeps its own selection at A12 but is not the active tab.
# ---------------------------------------------------------------------

$ws2.Select()
$ws2.Range("A12").Select()
$ws1.Select()
$ws1.Range("K4").Select()
$excel.ActiveWindow.Zoom = 100
